$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Arts and Humanities
$ws.Range("C2").Value = 9
$ws.Range("D2").Value = 0

# Row 5 - Council on Environmental Quality
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 0

# Row 6 - Department of Commerce
$ws.Range("C6").Value = 28
$ws.Range("D6").Value = 0

# Row 8 - Department of Education
$ws.Range("C8").Value = 25
$ws.Range("D8").Value = 0

# Row 9 - Department of Energy
$ws.Range("B9").Value = 22
$ws.Range("D9").Value = 1

# Row 10 - Department of Health and Human Services
$ws.Range("C10").Value = 48
$ws.Range("D10").Value = 0

# Row 12 - Department of Housing and Urban Development
$ws.Range("B12").Value = 23
$ws.Range("D12").Value = 2

# Row 13 - Department of Justice
$ws.Range("B13").Value = 33
$ws.Range("D13").Value = 1

# Row 16 - Department of the Interior
$ws.Range("B16").Value = 14
$ws.Range("D16").Value = 1

# Row 17 - Department of the Treasury
$ws.Range("B17").Value = 30
$ws.Range("D17").Value = 1

# Row 18 - Department of Transportation
$ws.Range("B18").Value = 22
$ws.Range("D18").Value = 1

# Row 19 - Department of Veterans Affairs
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 0

# Row 20 - Environmental Protection Agency
$ws.Range("B20").Value = 14
$ws.Range("D20").Value = 1

# Row 22 - Federal Communications Commission
$ws.Range("C22").Value = 5
$ws.Range("D22").Value = 0

# Row 32 - Office of Personnel Management
$ws.Range("B32").Value = 12
$ws.Range("D32").Value = 1

# Row 38 - United States Department of Agriculture
$ws.Range("B38").Value = 22
$ws.Range("C38").Value = 21
